$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '50.949.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.909.28'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '370.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.543'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.78%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -5.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.33%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('E13').Value = '  -4.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.362.98'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.34'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.898.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.921'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '50.898.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('E19').Value = '  -7.42%  '
$ws.Range('E20').Value = '  -4.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0941'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '258.24'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.28%  '
$ws.Range('E26').Value = '  -2.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.167'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.78%  '
$ws.Range('E30').Value = '  -5.05%  '
$ws.Range('E31').Value = '  -5.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('E33').Value = '  -4.37%  '
$ws.Range('E34').Value = '  -2.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.27'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '34.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.02%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0419'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.11%  '
$ws.Range('E40').Value = '  -4.98%  '
$ws.Range('E41').Value = '  -4.67%  '
$ws.Range('E42').Value = '  -6.68%  '
$ws.Range('E43').Value = '  -3.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '118.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.40%  '
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.012.51'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.79%  '
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.13'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.189.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.96%  '
$ws.Range('E51').Value = '  -1.02%  '
